$d = $word.ActiveDocument

# Locate the table row that discusses replacing PowerPoint lectures with VR,
# then grab its third cell ("Negative impact or issue" column), which is
# currently empty.
$found = $d.Content.Duplicate
$ok = $found.Find.Execute("suggested we use VR instead", $true, $false, $false, `
                           $false, $false, $true, 1, $false, "", 0)

$sourceCell = $found.Cells.Item(1)
$row = $sourceCell.Row
$targetCell = $row.Cells.Item(3)

# Insert the new sentence just before the cell's trailing paragraph mark.
$cellRange = $targetCell.Range
$insertPos = $cellRange.End - 1
$cellRange.Collapse(0)
$cellRange.Text = "VR is Expensive."

# Match the surrounding table formatting (Tahoma, 10pt). Re-fetch the range
# fresh for each property so the formatting calls don't clobber each other.
$formatRange1 = $d.Range($insertPos, $insertPos + 16)
$formatRange1.Font.Name = "Tahoma"

$formatRange2 = $d.Range($insertPos, $insertPos + 16)
$formatRange2.Font.Size = 10
